# Update the NBR and BAR sheets: shift the "Cutoff" column values by +4,
# update the "Reaction_number" column, and trim the trailing rows (17-20)
# so each sheet ends at row 16 instead of row 20.

$wb = $excel.ActiveWorkbook

$nbr = @(
    @(5, 96),
    @(6, 97),
    @(7, 98),
    @(8, 99),
    @(9, 97),
    @(10, 97),
    @(11, 98),
    @(12, 98),
    @(13, 98),
    @(14, 95),
    @(15, 94),
    @(16, 94),
    @(17, 94),
    @(18, 94),
    @(19, 93)
)

$bar = @(
    @(5, 578),
    @(6, 575),
    @(7, 574),
    @(8, 571),
    @(9, 571),
    @(10, 571),
    @(11, 567),
    @(12, 566),
    @(13, 566),
    @(14, 569),
    @(15, 566),
    @(16, 567),
    @(17, 565),
    @(18, 563),
    @(19, 566)
)

$ws1 = $wb.Worksheets.Item("NBR")
$ws2 = $wb.Worksheets.Item("BAR")

for ($i = 0; $i -lt $nbr.Count; $i++) {
    $r = $i + 2
    $ws1.Cells.Item($r, 2).Value = $nbr[$i][0]
    $ws1.Cells.Item($r, 3).Value = $nbr[$i][1]
}

for ($i = 0; $i -lt $bar.Count; $i++) {
    $r = $i + 2
    $ws2.Cells.Item($r, 2).Value = $bar[$i][0]
    $ws2.Cells.Item($r, 3).Value = $bar[$i][1]
}

# Remove the now-unused trailing rows (17-20) from both sheets.
$ws1.Rows("17:20").Delete()
$ws2.Rows("17:20").Delete()
